$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 17
$ws.Range("B5").Value = 45
$ws.Range("B6").Value = 120
$ws.Range("B7").Value = 14
